$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new 2013 publication row using the ORIGINAL (pre-insert)
# column layout: A=ID, B=Authors, C=Year, D=Title, E=Publication,
# F=PublicationLocation, G=WebLink. This mirrors the order the shared
# strings were appended in the source edit.
$ws.Range("A3").Value = "Lenz_et_al_2013"
$ws.Range("D3").Value = "Virologic response and characterisation of HCV genotype 2-6 in patients receiving TMC435 monotherapy (study TMC435-C202)"
$ws.Range("E3").Value = "Journal of Hepatology"
$ws.Range("F3").Value = "58(3):445-51"
$ws.Range("G3").Value = "http://www.sciencedirect.com/science/article/pii/S016882781200829X"
$ws.Range("B3").Value = "Lenz O, Vijgen L, Berke JM, Cummings MD, Fevery B, Peeters M, De Smedt G, Moreno C, Picchio G"
$ws.Range("C3").Value = 2013

# --- Insert a new "DisplayName" column right after the ID column, then
# populate the header and the two display-name values.
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "DisplayName"
$ws.Range("B2").Value = "Lenz et al., 2010"
$ws.Range("B3").Value = "Lenz et al., 2013"

# --- Match the widened A:B columns (DisplayName needs more room) while
# leaving the Authors/WebLink columns (now C / G) at their old widths.
$ws.Range("A1:B1").ColumnWidth = 20.83

# --- Cosmetic leftover formatting on the PublicationLocation cell of the
# new row (mirrors a pasted-hyperlink style remnant) plus a couple of
# stray formatted-but-empty cells one row below the data.
$ws.Range("G3").Font.Name = "Arial"
$ws.Range("G3").Font.Size = 11
$ws.Range("G3").Font.Color = 0

$ws.Range("G4").Font.Name = "Arial"
$ws.Range("G4").Font.Size = 11
$ws.Range("G4").Font.Color = 0

$ws.Range("C4").Font.Name = "Calibri"
$ws.Range("C4").Font.Size = 12

# --- Page setup (portrait, standard paper) and an updated selection.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B3").Select()
